$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data set (rows 80..184, columns D..R) shifts down by one row,
# a brand new observation is inserted at row 80, and the former last row
# (184) is pushed down to a brand new row 185 (with A:C copied across).

# 1) Read the existing block D80:R184 (105 rows x 15 cols) before touching anything.
$block = $ws.Range("D80:R184").Value2()

# 2) Write that block back shifted down by one row -> D81:R185
$ws.Range("D81:R185").Value2 = $block

# 3) Copy the A:C identifying columns down into the newly created row 185
#    (same market/region as the rest of the sheet).
$abc = $ws.Range("A184:C184").Value2()
$ws.Range("A185:C185").Value2 = $abc

# The brand-new row 185 has no formatting yet; give its Fecha cell (D185) the
# same date number format used by every other row's D column.
$ws.Range("D185").NumberFormat = $ws.Range("D184").NumberFormat

# 4) Overwrite row 80 with the new observation introduced at the top of the shift.
$ws.Range("D80").Value2 = 44467
$ws.Range("E80").Value2 = 4
$ws.Range("F80").Value2 = 100112032
$ws.Range("G80").Value2 = "Zapallo italiano"
$ws.Range("H80").Value2 = "Sin especificar"
$ws.Range("I80").Value2 = "Primera"
$ws.Range("J80").Value2 = 520
$ws.Range("K80").Value2 = 14000
$ws.Range("L80").Value2 = 15000
$ws.Range("M80").Value2 = 14500
$ws.Range("N80").Value2 = "$/caja 70 unidades"
$ws.Range("O80").Value2 = "Provincia de Limarí"
$ws.Range("P80").Value2 = 207
$ws.Range("Q80").Value2 = 70
$ws.Range("R80").Value2 = "Hortaliza"
